# Generate Report for Handoff
# Updates the localization-status report after a new handoff round:
#  - Overview sheet: bump "Latest HO Xliff Generate Date" for the files
#    that were just handed off (rows 4-7 -> 38f4354c.., 961c28b4.., 99f796b7.., 9e792af6..)
#  - zh-cn sheet: those same rows move from priority "low" to "ht" and get a
#    fresh "Latest Handoff Datetime"
#  - de-de sheet: those same rows also move from priority "low" to "ht"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-22 20:30:54"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-22 20:30:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-22 20:30:54"
